$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.975606546170582
$ws.Range("C2").Value = 0.2864795845893013
$ws.Range("D2").Value = 0.008536791027573543
$ws.Range("E2").Value = 0.04218811757406282
$ws.Range("F2").Value = 4.77886095270901
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("J2").Value = 0.1229205902999615
$ws.Range("L2").Value = 0.3265923476160779
$ws.Range("M2").Value = 0.7321708664111028
$ws.Range("N2").Value = 3.259197752174245
$ws.Range("B3").Value = 3.881884048324139
$ws.Range("C3").Value = 0.2634821436229231
$ws.Range("D3").Value = 0.007445546138448833
$ws.Range("E3").Value = 0.04177184720882821
$ws.Range("F3").Value = 4.757142986330976
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("J3").Value = 0.1231051025488905
$ws.Range("L3").Value = 0.325325740671893
$ws.Range("M3").Value = 0.7192830793392204
$ws.Range("N3").Value = 3.275525839364903
$ws.Range("B4").Value = 3.826681723177046
$ws.Range("C4").Value = 0.2495364413009611
$ws.Range("D4").Value = 0.006775747241597685
$ws.Range("E4").Value = 0.04151075920721059
$ws.Range("F4").Value = 4.746015564320658
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("J4").Value = 0.123223871941113
$ws.Range("L4").Value = 0.3246897639079336
$ws.Range("M4").Value = 0.7117789697486501
$ws.Range("N4").Value = 3.286377293792881
$ws.Range("B5").Value = 3.804774224173968
$ws.Range("C5").Value = 0.243896892001942
$ws.Range("D5").Value = 0.006502784756978741
$ws.Range("E5").Value = 0.04140296628153717
$ws.Range("F5").Value = 4.742034781527224
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("J5").Value = 0.1232736517028687
$ws.Range("L5").Value = 0.3244662581943274
$ws.Range("M5").Value = 0.7088236774715142
$ws.Range("N5").Value = 3.291006579176511
$ws.Range("B6").Value = 3.80117196715338
$ws.Range("C6").Value = 0.2429630575609849
$ws.Range("D6").Value = 0.006457456608661261
$ws.Range("E6").Value = 0.04138498249892208
$ws.Range("F6").Value = 4.741407180372192
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("J6").Value = 0.1232820010401992
$ws.Range("L6").Value = 0.3244312997042442
$ws.Range("M6").Value = 0.7083391512283299
$ws.Range("N6").Value = 3.291787771275096
$ws.Range("B7").Value = 3.826383892514343
$ws.Range("C7").Value = 0.2494602091702234
$ws.Range("D7").Value = 0.006772066125005693
$ws.Range("E7").Value = 0.04150931115203882
$ws.Range("F7").Value = 4.745959637885662
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("J7").Value = 0.1232245376962005
$ws.Range("L7").Value = 0.3246866052155681
$ws.Range("M7").Value = 0.7117386980587455
$ws.Range("N7").Value = 3.286438887534246
$ws.Range("B8").Value = 3.942803682930617
$ws.Range("C8").Value = 0.2785134594909664
$ws.Range("D8").Value = 0.008160438470966369
$ws.Range("E8").Value = 0.04204571897037912
$ws.Range("F8").Value = 4.770913588035484
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("J8").Value = 0.1229830748068634
$ws.Range("L8").Value = 0.3261262180840561
$ws.Range("M8").Value = 0.7276421536213817
$ws.Range("N8").Value = 3.264655974072696
$ws.Range("B9").Value = 4.189782338302109
$ws.Range("C9").Value = 0.3369008100935673
$ws.Range("D9").Value = 0.01088821020021413
$ws.Range("E9").Value = 0.04305472707341274
$ws.Range("F9").Value = 4.837433824306828
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("J9").Value = 0.1225529035408757
$ws.Range("L9").Value = 0.3300734479397534
$ws.Range("M9").Value = 0.7620846647675137
$ws.Range("N9").Value = 3.228514076392116
$ws.Range("B10").Value = 4.382768886533199
$ws.Range("C10").Value = 0.3807021077580544
$ws.Range("D10").Value = 0.0129002702936134
$ws.Range("E10").Value = 0.04377095092723415
$ws.Range("F10").Value = 4.897136646087631
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("J10").Value = 0.1222630884284368
$ws.Range("L10").Value = 0.3336591854334614
$ws.Range("M10").Value = 0.7893929252098033
$ws.Range("N10").Value = 3.205997497605452
$ws.Range("B11").Value = 4.473100869906489
$ws.Range("C11").Value = 0.4008342669046101
$ws.Range("D11").Value = 0.01381839404302099
$ws.Range("E11").Value = 0.04409156101890321
$ws.Range("F11").Value = 4.926673709018758
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("J11").Value = 0.1221369000661001
$ws.Range("L11").Value = 0.3354395026590709
$ws.Range("M11").Value = 0.8022555433902596
$ws.Range("N11").Value = 3.196637113979094
$ws.Range("B12").Value = 4.507674978685202
$ws.Range("C12").Value = 0.4084882060028008
$ws.Range("D12").Value = 0.01416655371507147
$ws.Range("E12").Value = 0.0442122372067697
$ws.Range("F12").Value = 4.938202442106927
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("J12").Value = 0.1220899254591457
$ws.Range("L12").Value = 0.336135112240683
$ws.Range("M12").Value = 0.8071898374476945
$ws.Range("N12").Value = 3.193220043669839
$ws.Range("B13").Value = 4.500212457105818
$ws.Range("C13").Value = 0.4068384347850724
$ws.Range("D13").Value = 0.01409154840370519
$ws.Range("E13").Value = 0.04418627971910816
$ws.Range("F13").Value = 4.935704209360836
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("J13").Value = 0.1221000062838629
$ws.Range("L13").Value = 0.3359843466874253
$ws.Range("M13").Value = 0.8061243209282054
$ws.Range("N13").Value = 3.193950290483215
$ws.Range("B14").Value = 4.475937929408587
$ws.Range("C14").Value = 0.4014633497993714
$ws.Range("D14").Value = 0.01384702714618413
$ws.Range("E14").Value = 0.04410150369531252
$ws.Range("F14").Value = 4.927615286687569
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("J14").Value = 0.1221330192122259
$ws.Range("L14").Value = 0.3354963011682202
$ws.Range("M14").Value = 0.8026602167245471
$ws.Range("N14").Value = 3.19635343071964
$ws.Range("B15").Value = 4.461116976903384
$ws.Range("C15").Value = 0.3981749225540625
$ws.Range("D15").Value = 0.0136973165014922
$ws.Range("E15").Value = 0.04404948110785512
$ws.Range("F15").Value = 4.922705398717312
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("J15").Value = 0.1221533460277242
$ws.Range("L15").Value = 0.3352001514737282
$ws.Range("M15").Value = 0.8005466295437884
$ws.Range("N15").Value = 3.19784204619576
$ws.Range("B16").Value = 4.376916779492433
$ws.Range("C16").Value = 0.3793906347689813
$ws.Range("D16").Value = 0.01284033242341565
$ws.Range("E16").Value = 0.04374989498915927
$ws.Range("F16").Value = 4.895254328031314
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("J16").Value = 0.1222714486314329
$ws.Range("L16").Value = 0.3335458382545369
$ws.Range("M16").Value = 0.788561195721023
$ws.Range("N16").Value = 3.206627030179902
$ws.Range("B17").Value = 4.325914935094318
$ws.Range("C17").Value = 0.367920454380112
$ws.Range("D17").Value = 0.0123153809605725
$ws.Range("E17").Value = 0.0435647874237679
$ws.Range("F17").Value = 4.879024262666434
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("J17").Value = 0.1223453462207651
$ws.Range("L17").Value = 0.3325691665524459
$ws.Range("M17").Value = 0.7813213671413024
$ws.Range("N17").Value = 3.212242790549539
$ws.Range("B18").Value = 4.296819118230644
$ws.Range("C18").Value = 0.361342531786363
$ws.Range("D18").Value = 0.0120137040931283
$ws.Range("E18").Value = 0.04345782765807549
$ws.Range("F18").Value = 4.869912875786014
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("J18").Value = 0.1223883820969602
$ws.Range("L18").Value = 0.3320214494425926
$ws.Range("M18").Value = 0.7771986051275022
$ws.Range("N18").Value = 3.215555858055126
$ws.Range("B19").Value = 4.287008782698422
$ws.Range("C19").Value = 0.3591186778294571
$ws.Range("D19").Value = 0.01191160396502511
$ws.Range("E19").Value = 0.0434215281288397
$ws.Range("F19").Value = 4.866866287940411
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("J19").Value = 0.1224030447351998
$ws.Range("L19").Value = 0.3318384131764986
$ws.Range("M19").Value = 0.7758098115236507
$ws.Range("N19").Value = 3.216691848298851
$ws.Range("B20").Value = 4.331319409556215
$ws.Range("C20").Value = 0.3691394602715832
$ws.Range("D20").Value = 0.01237123531925732
$ws.Range("E20").Value = 0.04358454311684046
$ws.Range("F20").Value = 4.880728814435173
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("J20").Value = 0.1223374246592095
$ws.Range("L20").Value = 0.3326716820229478
$ws.Range("M20").Value = 0.7820877732636191
$ws.Range("N20").Value = 3.21163638561417
$ws.Range("B21").Value = 4.483057961541022
$ws.Range("C21").Value = 0.4030413146423939
$ws.Range("D21").Value = 0.01391883513482384
$ws.Range("E21").Value = 0.04412642420055235
$ws.Range("F21").Value = 4.929981857672459
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("J21").Value = 0.1221233005445188
$ws.Range("L21").Value = 0.3356390700652696
$ws.Range("M21").Value = 0.8036759828175732
$ws.Range("N21").Value = 3.195644104939291
$ws.Range("B22").Value = 4.584369978590246
$ws.Range("C22").Value = 0.4253751936075219
$ws.Range("D22").Value = 0.01493315219644842
$ws.Range("E22").Value = 0.04447631648999106
$ws.Range("F22").Value = 4.964175683777142
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("J22").Value = 0.121988079420023
$ws.Range("L22").Value = 0.3377034097486558
$ws.Range("M22").Value = 0.8181553663516041
$ws.Range("N22").Value = 3.185935700009523
$ws.Range("B23").Value = 4.530101204993741
$ws.Range("C23").Value = 0.4134387802268975
$ws.Range("D23").Value = 0.01439150286874025
$ws.Range("E23").Value = 0.04428995656822732
$ws.Range("F23").Value = 4.94574183283612
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("J23").Value = 0.1220598182157464
$ws.Range("L23").Value = 0.3365901982489277
$ws.Range("M23").Value = 0.8103934897451168
$ws.Range("N23").Value = 3.191049029965484
$ws.Range("B24").Value = 4.328875343406082
$ws.Range("C24").Value = 0.3685882965550604
$ws.Range("D24").Value = 0.01234598318074376
$ws.Range("E24").Value = 0.04357561324786019
$ws.Range("F24").Value = 4.879957503037019
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("J24").Value = 0.1223410042803068
$ws.Range("L24").Value = 0.3326252918428168
$ws.Range("M24").Value = 0.781741158002589
$ws.Range("N24").Value = 3.211910278156722
$ws.Range("B25").Value = 4.120953628334007
$ws.Range("C25").Value = 0.3209497907130583
$ws.Range("D25").Value = 0.0101492271352015
$ws.Range("E25").Value = 0.04278624686193311
$ws.Range("F25").Value = 4.817544464093316
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("J25").Value = 0.122664656284277
$ws.Range("L25").Value = 0.3288852361704073
$ws.Range("M25").Value = 0.7524166107988961
$ws.Range("N25").Value = 3.237584629830593
